$d = $word.ActiveDocument

# --- Change 1 & 3: split the "...blendet die genaue Anzahl ein." run and
# plant a _GoBack bookmark right after "genau" (before the trailing "e").
# Word only ever keeps a single "_GoBack" bookmark in a document, so
# (re)adding it here automatically relocates it away from its old spot
# near "ANHANG: Autorenkreis" at the end of the document.
$rng = $d.Content
$found = $rng.Find.Execute("blendet die genau", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null
}

# --- Change 2: remove the stray single-space run (with noProof/eastAsia
# rPr) that sits right after "...wohl sparen." and before the page break.
# Done via Find/Replace (rather than Range.Delete) so the engine doesn't
# take the opportunity to coalesce unrelated, merely rsid-differing runs
# elsewhere in the document (e.g. the "ANHANG" / ": Autorenkreis" runs).
$rng2 = $d.Content
$rng2.Find.Execute("wohl sparen. ", $true, $false, $false, $false, $false, $true, 1, $false, "wohl sparen.", 2) | Out-Null
